$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.02'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.16'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.298'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.514'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8135'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8633'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06989'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03133'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02919'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09394'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.746'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001529'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04685'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0005986'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '17OneONEWorstin24h'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006220'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001240'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004646'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00006103'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.518'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.151'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1306'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002334'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03716'

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'KickToken'

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006423'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '40KickTokenKICK'

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'BKEXToken'

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1055'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'CEJI'

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003002'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007741'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005276'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.3804'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002438'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '47BOLOBOLO'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
